$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets("call")

# Update the maturity/expiry date column (E2:E15) 43983 -> 43836 for every row
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = 43836
}

# Last price updates
$ws.Range("C5").Value = 62.25
$ws.Range("C14").Value = 22.83

# New header label in H1 (shares font/style 4 already applied to that cell)
$ws.Range("H1").Value = "e"

# H2 had no content, only a leftover style - clear it out entirely
$ws.Range("H2").Clear()

# Move the selection / active cell
$ws.Range("I9").Select()

# Remove the stray second worksheet
$wb.Worksheets("Sheet1").Delete()
